$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")
$ws.Range("A554").NumberFormat = "@"
$ws.Range("A554").Value = "50060437"
